$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for row 2 (BKT0)
$ws.Range("B2").Value = 12973319
$ws.Range("C2").Value = 277
$ws.Range("D2").Value = 215
$ws.Range("E2").Value = 61
$ws.Range("I2").Value = 10370416
$ws.Range("J2").Value = 2546642
$ws.Range("M2").Value = 56261
$ws.Range("N2").Value = 79.94
$ws.Range("O2").Value = 19.63
$ws.Range("R2").Value = 0.43
$ws.Range("S2").Value = 19.63
$ws.Range("U2").Value = 1128903

# Updated values for row 3 (BKT1)
$ws.Range("B3").Value = 4619368
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 64
$ws.Range("G3").Value = 8
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4186884
$ws.Range("L3").Value = 280699
$ws.Range("M3").Value = 151785
$ws.Range("N3").Value = 90.64
$ws.Range("Q3").Value = 6.08
$ws.Range("R3").Value = 3.29
$ws.Range("S3").Value = 6.08
$ws.Range("T3").Value = 6.08
$ws.Range("U3").Value = 231682

# Updated values for row 4 (BKT2)
$ws.Range("B4").Value = 1375085
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 34
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 1262414
$ws.Range("J4").Value = 9970
$ws.Range("L4").Value = 102701
$ws.Range("N4").Value = 91.81
$ws.Range("O4").Value = 0.73
$ws.Range("Q4").Value = 7.47
$ws.Range("S4").Value = 8.199999999999999
$ws.Range("T4").Value = 7.47
$ws.Range("U4").Value = 63254

# Remove row 5 (BKT3) entirely - data no longer present
$ws.Rows.Item(5).Delete()
